$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data (A: value, B: date-time serial number with existing B-column style)
$data = @(
    @(3006.67, 43879.75540255231),
    @(3006.67, 43879.75543708817),
    @(3006.67, 43879.75547123926),
    @(3006.67, 43879.75550549213)
)

$row = 6
foreach ($pair in $data) {
    $ws.Cells.Item($row, 1).Value = $pair[0]
    $ws.Cells.Item($row, 2).Value = $pair[1]
    # Match the existing date/time number format used by the rest of column B
    $ws.Cells.Item($row, 2).NumberFormat = $ws.Cells.Item(2, 2).NumberFormat
    $row++
}
